$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.651.52"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").Value = "2.441.10"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.33%  "
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "2.820.34"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "2.447.93"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("D18").Value = "44.581.12"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0767"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "128.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.35%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0290"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").Value = "1.955.38"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.73%  "
